# Refresh the cryptocurrency price/volume snapshot on Sheet1 (cryptos.xlsx).
# Price values in column D are free-text strings (e.g. "3.802.70" uses '.'
# as a thousands separator) - set NumberFormat to Text ("@") first for the
# cells whose new value would otherwise be auto-parsed as a plain number,
# so the literal text is preserved exactly as scraped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.696.62'
$ws.Range("E2").Value = '  -2.66%  '
$ws.Range("D3").Value = '3.809.96'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.12'
$ws.Range("E5").Value = '  -2.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.36'
$ws.Range("E6").Value = '  -4.37%  '
$ws.Range("D7").Value = '3.807.03'
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("E10").Value = '  -4.32%  '
$ws.Range("E11").Value = '  -5.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").Value = '  -3.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.85'
$ws.Range("E13").Value = '  -2.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000245'
$ws.Range("E14").Value = '  -3.78%  '
$ws.Range("D15").Value = '4.442.64'
$ws.Range("E15").Value = '  +0.64%  '
$ws.Range("D16").Value = '3.804.97'
$ws.Range("E16").Value = '  +0.75%  '
$ws.Range("D17").Value = '67.773.92'
$ws.Range("E17").Value = '  -2.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.27'
$ws.Range("E18").Value = '  -4.22%  '
$ws.Range("E19").Value = '  -3.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.38'
$ws.Range("E20").Value = '  +5.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '495.80'
$ws.Range("E21").Value = '  -3.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.67'
$ws.Range("E22").Value = '  +2.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.745'
$ws.Range("E23").Value = '  +1.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.01'
$ws.Range("E24").Value = '  -0.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("E25").Value = '  -4.18%  '
$ws.Range("E26").Value = '  +6.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.41'
$ws.Range("E27").Value = '  -4.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.20'
$ws.Range("E28").Value = '  -3.66%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.45'
$ws.Range("E31").Value = '  -3.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.02'
$ws.Range("E32").Value = '  +6.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.93'
$ws.Range("E33").Value = '  -2.30%  '
$ws.Range("E34").Value = '  -4.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  -3.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.86'
$ws.Range("E37").Value = '  -4.70%  '
$ws.Range("E38").Value = '  -5.74%  '
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.330'
$ws.Range("E39").Value = '  -3.21%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '458.97'
$ws.Range("E40").Value = '  -0.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '49.16'
$ws.Range("E41").Value = '  -1.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.01'
$ws.Range("E42").Value = '  -3.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.86'
$ws.Range("E43").Value = '  -3.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.46'
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.62'
$ws.Range("E46").Value = '  -8.49%  '
$ws.Range("D47").Value = '2.849.54'
$ws.Range("E47").Value = '  -3.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '140.24'
$ws.Range("E48").Value = '  +0.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0354'
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.69'
$ws.Range("E50").Value = '  +14.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.93'
$ws.Range("E51").Value = '  -5.15%  '
